$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "Kesehatan Advice - Positive" to create the new
#        "Kesehatan Advice - Negative" sheet right after it, and make
#        it the active tab (mirrors tabSelected moving to the new
#        sheet / activeTab becoming 2). ---
$wsPositive = $wb.Worksheets.Item("Kesehatan Advice - Positive")
$wsPositive.Copy([System.Reflection.Missing]::Value, $wsPositive)
$wsNegative = $wb.Worksheets.Item(3)
$wsNegative.Name = "Kesehatan Advice - Negative"
$wsNegative.Select()

# --- 2. Populate the "Global" sheet with the automation-config header
#        row and its values (PASSWORD / COVER_TITLE / COVER_SUBTITLE /
#        COMPANY_LOGO / PROJECT_NAME / SCREEN_SHOOT). Values are written
#        in this particular order so the shared-string table comes out
#        in the same sequence as the authored workbook. ---
$wsGlobal = $wb.Worksheets.Item("Global")

$wsGlobal.Range("E1").Value = "PROJECT_NAME"
$wsGlobal.Range("B2").Value = "SuperApps-BSI Super Apps versi 1.0.2 (7761)"
$wsGlobal.Range("D1").Value = "COMPANY_LOGO"
$wsGlobal.Range("F1").Value = "SCREEN_SHOOT"
$wsGlobal.Range("B1").Value = "COVER_TITLE"
$wsGlobal.Range("D2").Value = "D:\Mentahan\PlugIn UFT\Napalm\Napalm\libray\BSILOGO.jpeg"
$wsGlobal.Range("E2").Value = "SuperApps-BSIMobile"
$wsGlobal.Range("A1").Value = "PASSWORD"
$wsGlobal.Range("C1").Value = "COVER_SUBTITLE"
$wsGlobal.Range("C2").Value = "Automation Testing - SuperApps-BSIMobile"

# Row 2 gets a top+bottom rule across the whole record, with an extra
# right-hand rule closing off the last column (F2) so the block reads
# like a bordered table row.
$rowRule = $wsGlobal.Range("A2:F2")
$topEdge = $rowRule.Borders.Item(8)
$topEdge.Color = 0
$topEdge.LineStyle = 1
$bottomEdge = $rowRule.Borders.Item(9)
$bottomEdge.Color = 0
$bottomEdge.LineStyle = 1

$rightEdge = $wsGlobal.Range("F2").Borders.Item(10)
$rightEdge.Color = 0
$rightEdge.LineStyle = 1

# Size the columns that now hold data to fit their contents.
$wsGlobal.Range("A1").ColumnWidth = 10
$wsGlobal.Range("B1").ColumnWidth = 26.666666666666668
$wsGlobal.Range("C1").ColumnWidth = 14.5
$wsGlobal.Range("D1").ColumnWidth = 51.5
$wsGlobal.Range("E1").ColumnWidth = 13.666666666666666
$wsGlobal.Range("F1").ColumnWidth = 13.5

Write-Output "done"
